# Auto-generated edit script: strip footnote markers like ' [1]' and replace
# embedded newlines with a single space across the vaccine price-list cells.
# Also removes the duplicate 'Afluria\nQuadrivalent' shared string so it collapses
# onto the existing 'Afluria Quadrivalent' entry (sheet 'Adult Influenza Vaccine').

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$updates1 = @{
    'A2' = 'Dengue Tetravalent Vaccine, Live '
    'A3' = 'DTaP '
    'A4' = 'DTaP '
    'A5' = 'DTaP-IPV '
    'A6' = 'DTaP-IPV '
    'A7' = 'DTaP-IPV '
    'A8' = 'DTaP-Hep B-IPV '
    'A9' = 'DTaP-IP-HI '
    'A10' = 'DTaP-IPV-HIB-HEPB '
    'A11' = 'DTaP-IPV-HIB-HEPB '
    'A12' = 'e-IPV '
    'A13' = 'Hepatitis A Pediatric '
    'A14' = 'Hepatitis A Pediatric '
    'A15' = 'Hepatitis A-Hepatitis B 18 only '
    'A16' = 'Hepatitis B  Pediatric/Adolescent'
    'A17' = 'Hepatitis B  Pediatric/Adolescent'
    'A18' = 'Hepatitis B  Pediatric/Adolescent'
    'A19' = 'Hib '
    'A20' = 'Hib '
    'A21' = 'Hib '
    'A22' = 'HPV - Human Papillomavirus 9-valent '
    'A23' = 'MENB - Meningococcal Group B '
    'A24' = 'MENB - Meningococcal Group B '
    'A25' = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
    'A26' = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
    'A27' = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
    'A28' = 'Measles, Mumps and Rubella (MMR) '
    'A29' = 'MMR/Varicella '
    'A30' = 'Pneumococcal 13-valent  (Pediatric)'
    'A32' = 'Rotavirus, Live, Oral, Pentavalent '
    'A33' = 'Rotavirus, Live, Oral, Pentavalent '
    'A34' = 'Rotavirus, Live, Oral, Oral '
    'A35' = 'Tetanus and Diphtheria Toxoids '
    'A36' = 'Tetanus and Diphtheria Toxoids '
    'A37' = 'Tetanus and Diphtheria Toxoids '
    'A38' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A39' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A40' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A41' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A42' = 'Varicella '
}
foreach ($addr in $updates1.Keys) {
    $ws1.Range($addr).Value = $updates1[$addr]
}

$ws2 = $wb.Worksheets.Item(2)
$updates2 = @{
    'A2' = 'Hepatitis A Adult '
    'A3' = 'Hepatitis A Adult '
    'A4' = 'Hepatitis A-Hepatitis B Adult '
    'A5' = 'Hepatitis B Adult '
    'A6' = 'Hepatitis B Adult '
    'A7' = 'Hepatitis B Adult '
    'A8' = 'HPV-Human Papillomavirus 9 Valent '
    'A9' = 'Measles, Mumps,  Rubella '
    'A10' = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
    'A11' = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
    'A12' = 'MENB - Meningococcal Group B '
    'A13' = 'MENB - Meningococcal Group B '
    'A14' = 'Pneumococcal 15-valent '
    'A15' = 'Pneumococcal 20-valent '
    'A17' = 'Tetanus and Diphtheria Toxoids '
    'A18' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A19' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A20' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A21' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A22' = 'Varicella '
}
foreach ($addr in $updates2.Keys) {
    $ws2.Range($addr).Value = $updates2[$addr]
}

$ws3 = $wb.Worksheets.Item(3)
$updates3 = @{
    'A2' = 'Influenza  (Age 6 months and older)'
    'B2' = 'Fluzone Quadrivalent'
    'A3' = 'Influenza  (Age 6 months and older)'
    'B3' = 'Fluzone Quadrivalent'
    'A4' = 'Influenza  (Age 6 months and older)'
    'B4' = 'Fluzone Quadrivalent'
    'A5' = 'Influenza  (Age 6 months and older)'
    'B5' = 'Fluarix Quadrivalent'
    'A6' = 'Influenza  (Age 6 months and older)'
    'B6' = 'FluLaval Quadrivalent'
    'A7' = 'Influenza  (Age 6 months and older)'
    'A8' = 'Influenza  (Age 6 months and older)'
    'A9' = 'Influenza  (Age 36 months and older)'
    'A10' = 'Influenza  (Age 6 months and older)'
    'A11' = 'Influenza  Live, Intranasal (Age 2-49 years)'
    'B11' = 'FluMist Quadrivalent'
}
foreach ($addr in $updates3.Keys) {
    $ws3.Range($addr).Value = $updates3[$addr]
}

$ws4 = $wb.Worksheets.Item(4)
$updates4 = @{
    'A2' = 'Influenza  (Age 6 months and older)'
    'B2' = 'Fluzone Quadrivalent'
    'A3' = 'Influenza  (Age 6 months and older)'
    'B3' = 'Fluzone Quadrivalent'
    'A4' = 'Influenza  (Age 6 months and older)'
    'B4' = 'Fluzone Quadrivalent'
    'A5' = 'Influenza  (Age 6 months and older)'
    'B5' = 'Fluarix Quadrivalent'
    'A6' = 'Influenza  (Age 6 months and older)'
    'B6' = 'FluLaval Quadrivalent'
    'A7' = 'Influenza  (Age 6 months and older)'
    'A8' = 'Influenza  (Age 6 months and older)'
    'A9' = 'Influenza  (Age 36 months and older)'
    'B9' = 'Afluria Quadrivalent'
    'A10' = 'Influenza  (Age 6 months and older)'
    'B10' = 'Afluria Quadrivalent'
    'A11' = 'Influenza  Live, Intranasal (Age 2-49 years)'
    'B11' = 'FluMist Quadrivalent'
}
foreach ($addr in $updates4.Keys) {
    $ws4.Range($addr).Value = $updates4[$addr]
}
